$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for rows 2-5 (A:E,G = text, F = number, H = boolean) ---
$ws.Range("A2").Value = "Squelette_sujet_Test_1"
$ws.Range("B2").Value = "Le"
$ws.Range("C2").Value = "melon"
$ws.Range("D2").Value = "Des"
$ws.Range("E2").Value = "melons"
$ws.Range("F2").Value = 1.0783318998292089
$ws.Range("G2").Value = "e"
$ws.Range("H2").Value = $false

$ws.Range("A3").Value = "Squelette_sujet_Test_2"
$ws.Range("B3").Value = "Les"
$ws.Range("C3").Value = "oignons"
$ws.Range("D3").Value = "Des"
$ws.Range("E3").Value = "oignons"
$ws.Range("F3").Value = 0.22487019980326295
$ws.Range("G3").Value = "q"
$ws.Range("H3").Value = $true

$ws.Range("A4").Value = "Squelette_sujet_Test_3"
$ws.Range("B4").Value = "La"
$ws.Range("C4").Value = "patate"
$ws.Range("D4").Value = "Une"
$ws.Range("E4").Value = "patate"
$ws.Range("F4").Value = 0.31397949997335672
$ws.Range("G4").Value = "q"
$ws.Range("H4").Value = $true

$ws.Range("A5").Value = "Squelette_sujet_Test_4"
$ws.Range("B5").Value = "La"
$ws.Range("C5").Value = "mangue"
$ws.Range("D5").Value = "Des"
$ws.Range("E5").Value = "mangues"
$ws.Range("F5").Value = 0.52217360027134418
$ws.Range("G5").Value = "e"
$ws.Range("H5").Value = $false

# --- I2: becomes a plain number 0 (was the text "ERREUR") while keeping its
# existing "text" cell style. Reset format to General first so the literal
# value is stored as a real number, then restore the Text format. ---
$ws.Range("I2").NumberFormat = "General"
$ws.Range("I2").Value = 0
$ws.Range("I2").NumberFormat = "@"

# --- Re-apply the "text" style to the header row and the text columns of
# rows 2-5 so they pick up a distinct style entry (border touch forces a
# fresh style record, matching the edit that produced the workbook). ---
$textRanges = @("A1:I1", "A2:E2", "G2", "A3:E3", "G3", "A4:E4", "G4", "A5:E5", "G5")
foreach ($rng in $textRanges) {
    $ws.Range($rng).Borders.LineStyle = 1
    $ws.Range($rng).NumberFormat = "@"
}

# --- Column width tweaks ---
$ws.Columns.Item(1).ColumnWidth = 22
$ws.Columns.Item(6).ColumnWidth = 12.7109375
$ws.Columns.Item(9).ColumnWidth = 6.5703125

Write-Output "Done"
